# Add a new worksheet named "iPhone-6 Raj" to the workbook, positioned after
# the last existing sheet ("null"). The new sheet reproduces the first two
# rows of the "null" sheet (same shared-string text + the green "PASS" fill
# used for passing test results), matching what a "duplicate sheet, trim to
# the Raj/iPhone-6 run, rename" edit would produce.

$wb = $excel.ActiveWorkbook

$sourceSheet = $wb.Worksheets.Item("null")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# Copy the "null" worksheet (brings along cell values + the PASS/FAIL fill
# styles) and drop it after the last sheet in the workbook.
$sourceSheet.Copy($null, $lastSheet)

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "iPhone-6 Raj"

# Keep only the first two rows (Test Parameters / Avner 3.01 header row and
# the SATeam PASS row) - drop the Avner/Sasha rows that came along with the
# copy.
$newSheet.Rows.Item(4).Delete()
$newSheet.Rows.Item(3).Delete()
